$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J, matching style of existing headers (s=1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (col 9) and J (col 10), rows 2-38
$data = @(
    @{Row=2; I=6; J=6},
    @{Row=3; I=9; J=9},
    @{Row=4; I=9; J=9},
    @{Row=5; I=9; J=9},
    @{Row=6; I=9; J=9},
    @{Row=7; I=9; J=9},
    @{Row=8; I=7; J=7},
    @{Row=9; I=7; J=7},
    @{Row=10; I=5; J=5},
    @{Row=11; I=9; J=9},
    @{Row=12; I=4; J=4},
    @{Row=13; I=9; J=9},
    @{Row=14; I=10; J=10},
    @{Row=15; I=8; J=8},
    @{Row=16; I=7; J=8},
    @{Row=17; I=7; J=7},
    @{Row=18; I=7; J=7},
    @{Row=19; I=9; J=9},
    @{Row=20; I=7; J=8},
    @{Row=21; I=8; J=9},
    @{Row=22; I=7; J=7},
    @{Row=23; I=7; J=7},
    @{Row=24; I=8; J=8},
    @{Row=25; I=1; J=3},
    @{Row=26; I=9; J=9},
    @{Row=27; I=9; J=9},
    @{Row=28; I=7; J=7},
    @{Row=29; I=8; J=8},
    @{Row=30; I=6; J=6},
    @{Row=31; I=9; J=9},
    @{Row=32; I=4; J=4},
    @{Row=33; I=9; J=9},
    @{Row=34; I=5; J=6},
    @{Row=35; I=4; J=4},
    @{Row=36; I=4; J=4},
    @{Row=37; I=8; J=8},
    @{Row=38; I=4; J=4}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 9).Value = $entry.I
    $ws.Cells.Item($entry.Row, 10).Value = $entry.J
}
